$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 58 (this shifts rows 58+ down by one,
# which pushes the Staff table's trailing rows and every table below it down)
$ws.Rows("58:58").Insert()

# Copy formatting from row 57 (the original last "Staff" row) into new row 58
$ws.Range("G57:K57").Copy()
$ws.Range("G58:K58").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 58 now holds the content that used to be in row 57 (the "isDeleted" field)
$ws.Range("G58").Value2 = 13
$ws.Range("H58").Value2 = "isDeleted"
$ws.Range("I58").Value2 = "bit"
$ws.Range("J58").Value2 = "NOT NULL, DEFAULT = 0"
$ws.Range("K58").Value2 = "Trạng thái xóa, mặc định là false"

# Row 57 becomes the new "isActived" field describing staff activation status
$ws.Range("H57").Value2 = "isActived"
$ws.Range("K57").Value2 = "Trạng thái kích hoạt tài khoản, 0 = chưa, 1 = rồi"

# Update the view: scroll position and active selection
[void]$ws.Select()
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 39
$win.ScrollColumn = 4
[void]$ws.Range("K57").Select()

Write-Host "done"
